$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.7146045689148273
$ws.Range("B1").Value = 3.5604435535863543
$ws.Range("C1").Value = 133480

$ws.Range("A4").Value = 2.6850329762824661
$ws.Range("B4").Value = 2.8815188429407783
$ws.Range("C4").Value = 200615
